$wb = $excel.ActiveWorkbook

# --- Sheet "readme": reorder the JobNo/Date/Author/sheet_name columns ---
$ws = $wb.Worksheets.Item("readme")

# Capture current values before overwriting anything (columns B..E for rows 1..12)
$numRows = 12
$authorCol = @()
$sheetNameCol = @()
$jobNoCol = @()
$dateCol = @()
for ($r = 1; $r -le $numRows; $r++) {
    $authorCol += , $ws.Cells.Item($r, 2).Value2
    $sheetNameCol += , $ws.Cells.Item($r, 3).Value2
    $jobNoCol += , $ws.Cells.Item($r, 4).Value2
    $dateCol += , $ws.Cells.Item($r, 5).Value2
}

# New layout: B=JobNo, C=Date, D=Author, E=sheet_name
for ($r = 1; $r -le $numRows; $r++) {
    $i = $r - 1
    $ws.Cells.Item($r, 2).Value = $jobNoCol[$i]
    $ws.Cells.Item($r, 3).Value = $dateCol[$i]
    $ws.Cells.Item($r, 4).Value = $authorCol[$i]
    $ws.Cells.Item($r, 5).Value = $sheetNameCol[$i]
}

# Update the Date values (rows 2-12) from 20220225 to 20220228
for ($r = 2; $r -le $numRows; $r++) {
    $ws.Cells.Item($r, 3).Value = "20220228"
}

# --- Sheet "Project Information": update Date of Analysis value ---
$wsInfo = $wb.Worksheets.Item("Project Information")
$wsInfo.Range("B11").Value = "2022-02-28 12:43:33.002303"
